{"js": "// The document's first paragraph holds a date stamp, followed by a single\n// 20-row x 5-column table of arithmetic problems (\"a+b=c\" / \"a-b=c\").\n// The commit regenerates the whole worksheet: the date advances one day\n// and every cell in the table gets a freshly generated problem. Crucially,\n// the new values are NOT a simple re-shuffle of the old ones (several new\n// values coincide with old values that live elsewhere in the table), so we\n// must replace every value by its *position* in the document rather than\n// by searching for old text and replacing it \u2014 a text-based find/replace\n// could clobber a cell that was already updated to a value another cell\n// used to hold.\n\nconst newDate = \"2025-01-02 Thursday\";\n\n// Row-major replacement values for the 20x5 table (100 cells total).\nconst newTableValues = [\n  [\"5+12=17\", \"58-0=58\", \"62+20=82\", \"78-75=3\", \"88+5=93\"],\n  [\"17+74=91\", \"97-29=68\", \"65-28=37\", \"27+44=71\", \"39+57=96\"],\n  [\"50-5=45\", \"28+39=67\", \"16+25=41\", \"3+91=94\", \"31-31=0\"],\n  [\"37+45=82\", \"10+32=42\", \"1+23=24\", \"65-49=16\", \"60-48=12\"],\n  [\"18+17=35\", \"42+17=59\", \"42+17=59\", \"91+3=94\", \"27+33=60\"],\n  [\"11+25=36\", \"92-50=42\", \"52+17=69\", \"43+15=58\", \"9+24=33\"],\n  [\"52-10=42\", \"48-22=26\", \"62+15=77\", \"5+35=40\", \"16+70=86\"],\n  [\"97-93=4\", \"17+77=94\", \"62-50=12\", \"35-0=35\", \"12+34=46\"],\n  [\"59-49=10\", \"58-43=15\", \"57+1=58\", \"58-51=7\", \"67+17=84\"],\n  [\"20+51=71\", \"97-24=73\", \"28+18=46\", \"77+15=92\", \"84-65=19\"],\n  [\"69-34=35\", \"29+65=94\", \"99-54=45\", \"41-35=6\", \"87-54=33\"],\n  [\"66-42=24\", \"37+1=38\", \"23+19=42\", \"37+13=50\", \"9+80=89\"],\n  [\"47+20=67\", \"28-9=19\", \"31+41=72\", \"44+13=57\", \"29+53=82\"],\n  [\"16+43=59\", \"46-23=23\", \"42-16=26\", \"43+37=80\", \"81-33=48\"],\n  [\"59-23=36\", \"42-27=15\", \"52+41=93\", \"11+83=94\", \"72-5=67\"],\n  [\"92-13=79\", \"80-37=43\", \"89-56=33\", \"23+53=76\", \"75-13=62\"],\n  [\"52-27=25\", \"9+73=82\", \"27+43=70\", \"71-39=32\", \"12+57=69\"],\n  [\"61-60=1\", \"5+45=50\", \"47-34=13\", \"89-83=6\", \"91-55=36\"],\n  [\"57+39=96\", \"27+26=53\", \"4+59=63\", \"94-61=33\", \"42+45=87\"],\n  [\"80-45=35\", \"51+41=92\", \"11+72=83\", \"75-71=4\", \"93+2=95\"],\n];\n\nconst body = context.document.body;\n\n// Update the date line (first paragraph in the body).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(newDate, Word.InsertLocation.replace);\n\n// Update every cell of the table in one shot, preserving per-cell\n// formatting (font/size/alignment) since only the text run content\n// changes, not the structure.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\ntables.items[0].values = newTableValues;\n\nawait context.sync();\n", "ps1": "# The document's first paragraph holds a date stamp, followed by a single\n# 20-row x 5-column table of arithmetic problems (\"a+b=c\" / \"a-b=c\").\n# The commit regenerates the whole worksheet: the date advances one day\n# and every cell in the table gets a freshly generated problem. Crucially,\n# the new values are NOT a simple re-shuffle of the old ones (several new\n# values coincide with old values that live elsewhere in the table), so we\n# must replace every value by its *position* (paragraph / row / column)\n# rather than by searching for old text and replacing it - a text-based\n# find/replace could clobber a cell that was already updated to a value\n# another cell used to hold.\n\n$d = $word.ActiveDocument\n\n# Update the date line (first paragraph in the body).\n$d.Paragraphs.Item(1).Range.Text = \"2025-01-02 Thursday\"\n\n# Row-major replacement values for the 20x5 table (100 cells total).\n$newTableValues = @(\n    @(\"5+12=17\",\"58-0=58\",\"62+20=82\",\"78-75=3\",\"88+5=93\"),\n    @(\"17+74=91\",\"97-29=68\",\"65-28=37\",\"27+44=71\",\"39+57=96\"),\n    @(\"50-5=45\",\"28+39=67\",\"16+25=41\",\"3+91=94\",\"31-31=0\"),\n    @(\"37+45=82\",\"10+32=42\",\"1+23=24\",\"65-49=16\",\"60-48=12\"),\n    @(\"18+17=35\",\"42+17=59\",\"42+17=59\",\"91+3=94\",\"27+33=60\"),\n    @(\"11+25=36\",\"92-50=42\",\"52+17=69\",\"43+15=58\",\"9+24=33\"),\n    @(\"52-10=42\",\"48-22=26\",\"62+15=77\",\"5+35=40\",\"16+70=86\"),\n    @(\"97-93=4\",\"17+77=94\",\"62-50=12\",\"35-0=35\",\"12+34=46\"),\n    @(\"59-49=10\",\"58-43=15\",\"57+1=58\",\"58-51=7\",\"67+17=84\"),\n    @(\"20+51=71\",\"97-24=73\",\"28+18=46\",\"77+15=92\",\"84-65=19\"),\n    @(\"69-34=35\",\"29+65=94\",\"99-54=45\",\"41-35=6\",\"87-54=33\"),\n    @(\"66-42=24\",\"37+1=38\",\"23+19=42\",\"37+13=50\",\"9+80=89\"),\n    @(\"47+20=67\",\"28-9=19\",\"31+41=72\",\"44+13=57\",\"29+53=82\"),\n    @(\"16+43=59\",\"46-23=23\",\"42-16=26\",\"43+37=80\",\"81-33=48\"),\n    @(\"59-23=36\",\"42-27=15\",\"52+41=93\",\"11+83=94\",\"72-5=67\"),\n    @(\"92-13=79\",\"80-37=43\",\"89-56=33\",\"23+53=76\",\"75-13=62\"),\n    @(\"52-27=25\",\"9+73=82\",\"27+43=70\",\"71-39=32\",\"12+57=69\"),\n    @(\"61-60=1\",\"5+45=50\",\"47-34=13\",\"89-83=6\",\"91-55=36\"),\n    @(\"57+39=96\",\"27+26=53\",\"4+59=63\",\"94-61=33\",\"42+45=87\"),\n    @(\"80-45=35\",\"51+41=92\",\"11+72=83\",\"75-71=4\",\"93+2=95\")\n)\n\n$t = $d.Tables.Item(1)\nfor ($r = 1; $r -le $newTableValues.Count; $r++) {\n    $rowValues = $newTableValues[$r - 1]\n    for ($c = 1; $c -le $rowValues.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
